# "update all reports + commit log"
# Appends 12 new log rows (54-65) to the commit-log table on Sheet1, two rows
# per work-date (Ori / Omri), then moves the sheet's view/selection down to
# the freshly-added tail of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- row data -------------------------------------------------------------
# NAME, DATE, FILE, COMMENT
$rows = @{
    54 = @("Omri", "14.4.19", "address_translation.v", "advance address protection coded")
    55 = @("Ori",  "14.4.19", "lru.v",       "LRU supports partial maping")
    56 = @("Ori",  "21.4.19", "m_handler.v", "MISS HANDLER support partial mapping")
    57 = @("Omri", "21.4.19", "all files",   "FIXing lintra problems")
    58 = @("Ori",  "28.4.19", "all cache",   "Mapping testbanch")
    59 = @("Omri", "28.4.19", "all AT",      "PROTECTION testbanch")
    60 = @("Ori",  "5.5.19",  "all files",   "code coverage and tools cleaning")
    61 = @("Omri", "5.5.19",  "all files",   "code coverage and tools cleaning")
    62 = @("Ori",  "12.5.19", "all files",   "establishing rules for Formal varification")
    63 = @("Omri", "12.5.19", "all files",   "establishing rules for Formal varification")
    64 = @("Ori",  "19.5.19", "all files",   "Debug formal varification fails")
    65 = @("Omri", "19.5.19", "all files",   "Debug formal varification fails")
}

# --- fill the DATE column first (one work-week's pair of rows at a time) --
# Force the "@" text format before writing so ambiguous D.M.YY strings (both
# parts <= 12, e.g. "5.5.19") aren't auto-parsed into date serials, then
# drop the format again so the cell is left with the workbook's default
# (unstyled) cell format, matching a plain text entry.
$dateRowOrder = @(54, 55, 56, 57, 58, 59, 64, 65, 62, 63, 60, 61)
foreach ($r in $dateRowOrder) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $rows[$r][1]
    $cell.ClearFormats()
}

# --- then fill NAME, COMMENT, FILE for every new row, top to bottom -------
foreach ($r in 54..65) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 4).Value = $vals[3]
    $ws.Cells.Item($r, 3).Value = $vals[2]
}

# --- move the view to the new bottom of the table --------------------------
$ws.Range("D66").Select()
